$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update E22 hours value from 11 to 15
$ws.Range("E22").Value = 15

# Update F22 activity text to append the new work description
$ws.Range("F22").Value = "Finalized sesssion closing. Bug fixes. Cleaned up unused code. Added some documentation. Migrated to a different Stored Procedure to close sessions, changed tests accordingly. Re-wrote inserting panic responses since Luca and Mihai obviously never tested it. Finally completed web sockets! Finished HTTPS. Hosted the Spring Boot server on the remote server. Fixed bugs and re-did the architecture notebook and updated all diagrams to reflect the release state of the application."

# Move the active cell selection from E23 to F23
$ws.Range("F23").Select()
